$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: first dish name changes from "суп из чечевицы с овощами" to "супец пупец"
# (set first so it reuses the shared-string slot vacated by the old B1 text)
$ws.Range("B2").Value = "супец пупец"

# B1: date range string changes from "23.04.2018-30.07.2018" to "30.04.2018-30.07.2018"
$ws.Range("B1").Value = "30.04.2018-30.07.2018"

# Update selection to B2
$ws.Range("B2").Select()
